$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point drift on the previous last row's timestamp
$ws.Range("A93").Value = 44406.76769951273

# Append the newly scraped row of data
$ws.Range("A94").Value = 44407.76788895095
$ws.Range("B94").Value = 80907
$ws.Range("C94").Value = 68298
$ws.Range("D94").Value = 3568
$ws.Range("E94").Value = 2269
$ws.Range("F94").Value = 1652
$ws.Range("G94").Value = 21299
$ws.Range("H94").Value = 1618
$ws.Range("I94").Value = 927
$ws.Range("J94").Value = 198
